$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# The deck currently has 2 slides:
#   1. "Christmas Presentation" (title slide)
#   2. "Christmas Tree" (title + picture)
#
# We need to insert a new "Winter Facts" slide between them, so the final
# order becomes:
#   1. "Christmas Presentation"
#   2. "Winter Facts"            (new)
#   3. "Christmas Tree"          (previous slide 2, pushed down)
# ---------------------------------------------------------------------------

$oldSlide2 = $p.Slides.Item(2)
$oldTitleShape = $oldSlide2.Shapes.Item(1)
$oldPicShape = $oldSlide2.Shapes.Item(2)

# Remember the old title's text/formatting so we can reproduce it on the
# slide the picture ends up on.
$oldTitleText = $oldTitleShape.TextFrame.TextRange.Text

# ---------------------------------------------------------------------------
# Step 1: capture a copy of the Christmas-tree picture before we start
# rearranging shapes (Copy/Paste keeps the embedded image relationship).
# ---------------------------------------------------------------------------
$oldPicShape.Copy()

# ---------------------------------------------------------------------------
# Step 2: append a brand-new slide using the "Title and Content" layout.
# This slide will become the new, final "Christmas Tree" slide (slide 3).
# ---------------------------------------------------------------------------
$treeSlide = $p.Slides.Add($p.Slides.Count + 1, 6)

$treeSlide.FollowMasterBackground = $false
$treeSlide.Background.Fill.ForeColor.RGB = 0x23380C
$treeSlide.Background.Fill.Solid()

$treeTitle = $treeSlide.Shapes.Item(1).TextFrame.TextRange
$treeTitle.Text = $oldTitleText
$treeTitle.ParagraphFormat.Alignment = 2
$treeTitle.Font.Size = 40
$treeTitle.Font.Bold = $true
$treeTitle.Font.Color.RGB = 0xFFFFFF
$treeTitle.Font.Name = "Calibri"

$treeSlide.Shapes.Paste()

# ---------------------------------------------------------------------------
# Step 3: turn the old slide 2 into the new "Winter Facts" slide. Clear its
# existing shapes (title + picture) and replace them with a clean
# Title + Content Placeholder pair (harvested from a scratch slide using the
# same layout), so no left-over formatting survives.
# ---------------------------------------------------------------------------
$oldSlide2.Shapes.Item(2).Delete()
$oldSlide2.Shapes.Item(1).Delete()
$oldSlide2.Shapes.Item(1).Delete()

$scratch = $p.Slides.Add($p.Slides.Count + 1, 2)
$scratchTitle = $scratch.Shapes.Item(1)
$scratchBody = $scratch.Shapes.Item(2)

$scratchTitle.Copy()
$oldSlide2.Shapes.Paste()
$scratchBody.Copy()
$oldSlide2.Shapes.Paste()

$scratch.Delete()

$winterTitle = $oldSlide2.Shapes.Item(1).TextFrame.TextRange
$winterTitle.Text = "Winter Facts"
$winterTitle.ParagraphFormat.Alignment = 1
$winterTitle.Font.Size = 40
$winterTitle.Font.Color.RGB = 0xFFFFFF
$winterTitle.Font.Name = "Calibri"

$winterBody = $oldSlide2.Shapes.Item(2).TextFrame.TextRange
$winterBody.Text = "1. Snow is white`r2. It's cold`r3. People build snowmen"
$winterBody.ParagraphFormat.Alignment = 1
$winterBody.Font.Size = 24
$winterBody.Font.Color.RGB = 0xFFFFFF
$winterBody.Font.Name = "Calibri"

Write-Output "Final slide count: $($p.Slides.Count)"
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slideTitle = $p.Slides.Item($i).Shapes.Item(1).TextFrame.TextRange.Text
    Write-Output "Slide $i -> $slideTitle"
}
